# Nexial showcase workbook update:
# Add a new "base" function outputToCloud(resource) and a new "text"
# function category containing spellCheck(var,profile,text) to the
# hidden '#system' reference sheet, and wire up the corresponding
# defined names.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------
# 1) Insert a brand-new column at Y (25th column). This shifts the
#    existing Y..AD columns (web, webalert, webcookie, ws, ws.async,
#    xml) one column to the right (Z..AE), matching the diff.
# ---------------------------------------------------------------
$ws.Columns.Item(25).Insert()

# New "text" function category lives in the freshly inserted column Y.
$ws.Range("Y1").Value = "text"
$ws.Range("Y2").Value = "spellCheck(var,profile,text)"

# ---------------------------------------------------------------
# 2) Insert "outputToCloud(resource)" into the "base" function list
#    (column E) in alphabetical order -- it sorts right before
#    "prependText(var,prependWith)" which currently sits at E22.
#    Shift only the column-E cells (E22:E38) down one row so the
#    neighbouring columns on those rows are left untouched.
# ---------------------------------------------------------------
for ($r = 38; $r -ge 22; $r--) {
    $src = $ws.Cells.Item($r, 5)
    $dst = $ws.Cells.Item($r + 1, 5)
    $dst.Value = $src.Value()
}
$ws.Cells.Item(22, 5).Value = "outputToCloud(resource)"

# ---------------------------------------------------------------
# 3) Insert "text" into the "target" list of category names
#    (column A) in alphabetical order -- it sorts right before
#    "web" which currently sits at A25. Shift only column-A cells
#    (A25:A30) down one row.
# ---------------------------------------------------------------
for ($r = 30; $r -ge 25; $r--) {
    $src = $ws.Cells.Item($r, 1)
    $dst = $ws.Cells.Item($r + 1, 1)
    $dst.Value = $src.Value()
}
$ws.Cells.Item(25, 1).Value = "text"

# ---------------------------------------------------------------
# 4) Fix up the defined names so every named range still points at
#    the right block of cells after the insertions above.
# ---------------------------------------------------------------
$wb.Names.Item("base").RefersTo = "='#system'!`$E`$2:`$E`$39"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$31"
$wb.Names.Item("web").RefersTo = "='#system'!`$Z`$2:`$Z`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AC`$2:`$AC`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AD`$2:`$AD`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AE`$2:`$AE`$27"

# Brand-new defined name for the "text" function category.
$wb.Names.Add("text", "='#system'!`$Y`$2:`$Y`$2")
